$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "261.59"
Set-TextValue "E2" "1.07%"
Set-TextValue "E3" "1.08%"
Set-TextValue "D4" "4.710"
Set-TextValue "E4" "0.66%"
Set-TextValue "D5" "0.06196"
Set-TextValue "E5" "3.23%"
Set-TextValue "D6" "6.714"
Set-TextValue "E6" "0.75%"
Set-TextValue "D7" "0.8507"
Set-TextValue "E7" "-0.80%"
Set-TextValue "D8" "0.9141"
Set-TextValue "E8" "-1.25%"
Set-TextValue "D9" "0.1411"
Set-TextValue "E9" "1.46%"
Set-TextValue "D10" "0.04515"
Set-TextValue "E10" "-5.96%"
Set-TextValue "D11" "0.07084"
Set-TextValue "E11" "1.02%"
Set-TextValue "D12" "0.03132"
Set-TextValue "E12" "0.28%"
Set-TextValue "D13" "0.09041"
Set-TextValue "D14" "0.001537"
Set-TextValue "E14" "0.31%"
Set-TextValue "D15" "0.0006165"
Set-TextValue "E15" "1.88%"
Set-TextValue "D16" "0.006055"
Set-TextValue "E16" "-0.40%"
Set-TextValue "D17" "3.460"
Set-TextValue "E17" "-0.02%"
Set-TextValue "D18" "3.165"
Set-TextValue "E18" "-0.01%"
Set-TextValue "E19" "1.33%"
Set-TextValue "D21" "0.1309"
Set-TextValue "E21" "0.86%"
Set-TextValue "D22" "4.083"
Set-TextValue "E22" "-0.95%"
Set-TextValue "E23" "-0.13%"
Set-TextValue "D24" "0.001216"
Set-TextValue "E24" "0.02%"
Set-TextValue "E25" "-5.71%"
Set-TextValue "E26" "0.09%"
Set-TextValue "D40" "0.03938"
Set-TextValue "E40" "2.57%"
Set-TextValue "E41" "-0.12%"
Set-TextValue "D42" "0.004131"
Set-TextValue "E42" "7.47%"
Set-TextValue "E43" "-9.70%"
Set-TextValue "D44" "0.01386"
Set-TextValue "E44" "-9.23%"
Set-TextValue "D45" "0.00005141"
Set-TextValue "E45" "0.60%"
Set-TextValue "E46" "0.10%"
Set-TextValue "D48" "0.1667"
Set-TextValue "E48" "10.88%"
Set-TextValue "E49" "0.10%"
Set-TextValue "E50" "0.10%"
